$p = $ppt.ActivePresentation

# Slide 2 (Introduction) - bullet wording tweaks
$s2 = $p.Slides.Item(2)
$s2body = $s2.Shapes.Item(2).TextFrame
$s2body.DeleteText()
$s2body.TextRange.Text = "- Addressing the problem of Generalized Category Discovery (GCD)`r- Introducing Expert-Contrastive Learning (XCon) for mining useful information from images`r- Utilizing k-means clustering and contrastive learning to learn discriminative features"

# Slide 3 (Problem Statement -> Motivation)
$s3 = $p.Slides.Item(3)
$s3.Shapes.Item(1).TextFrame.TextRange.Text = "Motivation"
$s3body = $s3.Shapes.Item(2).TextFrame
$s3body.DeleteText()
$s3body.TextRange.Text = "- Importance of Generalized Category Discovery in real-world applications`r- Contrasting unsupervised representation clustering with class-irrelevant cues vs. XCon's discriminative feature learning"

# Slide 4 (XCon Methodology -> XCon Method)
$s4 = $p.Slides.Item(4)
$s4.Shapes.Item(1).TextFrame.TextRange.Text = "XCon Method"
$s4body = $s4.Shapes.Item(2).TextFrame
$s4body.DeleteText()
$s4body.TextRange.Text = "- Expert-Contrastive Learning approach with k-means partitioning into expert sub-datasets`r- Eliminating negative influence of class-irrelevant cues for fine-grained category discovery`r- Using XCon for learning discriminative features and discovering new object categories"

# Slide 5 (Novel Category Discovery -> Contrastive Learning in GCD)
$s5 = $p.Slides.Item(5)
$s5.Shapes.Item(1).TextFrame.TextRange.Text = "Contrastive Learning in GCD"
$s5body = $s5.Shapes.Item(2).TextFrame
$s5body.DeleteText()
$s5body.TextRange.Text = "- Leveraging contrastive learning for effective representation learning in GCD`r- Using k-means grouping on self-supervised features for informative pairs`r- Focus on fine-grained category discovery with feature partitioning and contrastive pairs creation"

# Slide 6 (Previous Work -> Fine-grained GCD Challenges)
$s6 = $p.Slides.Item(6)
$s6.Shapes.Item(1).TextFrame.TextRange.Text = "Fine-grained GCD Challenges"
$s6body = $s6.Shapes.Item(2).TextFrame
$s6body.DeleteText()
$s6body.TextRange.Text = "- Detailed discriminative traits requirement for representation learning`r- Utilizing self-supervised representations and k-means clustering for feature integration`r- Applying supervised and self-supervised contrastive losses for feature refinement"

# Slide 7 (Experimental Results (Generic Datasets) -> Results on Generic Datasets)
$s7 = $p.Slides.Item(7)
$s7.Shapes.Item(1).TextFrame.TextRange.Text = "Results on Generic Datasets"
$s7body = $s7.Shapes.Item(2).TextFrame
$s7body.DeleteText()
$s7body.TextRange.Text = "- Evaluation on CIFAR10, CIFAR100, Stanford Cars, and ImageNet benchmarks`r- Performance improvement comparison with state-of-the-art methods`r- Application of semi-supervised k-means for evaluation metrics"

# Slide 8 (Experimental Results (Fine-grained Datasets) -> Results on Fine-grained Datasets)
$s8 = $p.Slides.Item(8)
$s8.Shapes.Item(1).TextFrame.TextRange.Text = "Results on Fine-grained Datasets"
$s8body = $s8.Shapes.Item(2).TextFrame
$s8body.DeleteText()
$s8body.TextRange.Text = "- Enhanced performance on CUB-200 and Standford Cars`r- Comparison of XCon with baseline under different parameters`r- Robust effectiveness analysis of XCon on fine-grained category discovery"

# Slide 9 (Qualitative Analysis) - bullet wording tweaks
$s9 = $p.Slides.Item(9)
$s9body = $s9.Shapes.Item(2).TextFrame
$s9body.DeleteText()
$s9body.TextRange.Text = "- Visualization of feature clustering in CIFAR10 using XCon`r- Clear boundaries between distinct categories with XCon's discriminative features`r- Improved categorization based on fine-grained features compared to DINO"

# Slide 10 (Conclusion) - bullet wording tweaks
$s10 = $p.Slides.Item(10)
$s10body = $s10.Shapes.Item(2).TextFrame
$s10body.DeleteText()
$s10body.TextRange.Text = "- Proposing XCon for Generalized Category Discovery with self-supervised representation`r- Learning fine-grained discriminative features for category discovery`r- Validation of XCon's effectiveness through performance improvements in experiments"
